$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Update the cached "datetimeFigureOut" field text on the slide master
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "7/1/2023"
    }
}

# Update the same cached field text on every slide layout
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "7/1/2023"
        }
    }
}

# Slide 16: fix "mis diarization" -> "mis-diarization" in the last bullet
$s = $p.Slides.Item(16)
$sh2 = $s.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$para5 = $tr2.Paragraphs(5, 1)
$para5.Text = "in analyzation calculate mis-diarization per second per segment then total"
